$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (old D = "Tipo") so "Tipo" shifts to column E
$ws.Columns("D").Insert()

# New column D header and value
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 0.1313832720735846

# Update existing MSE / R2 values
$ws.Range("B2").Value = 0.03272455184303397
$ws.Range("C2").Value = 0.9996627150218961
